$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# "ImportController" -> "StoringController": replace the "Import" prefix
# with "Storing", leaving "Controller" untouched. This naturally splits
# the single run into two runs ("Storing" / "Controller") the same way
# PowerPoint does when you retype part of a word.
[void]$tr.Replace("Import", "Storing")
